# Updates the 2022 full-year ERT_SU_CZ sheet with the final (actual) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ERT_SU_CZ")

# Release date (B2) moves from 2023-01-18 to 2023-04-18.
$ws.Range("B2").Value = 45034

# Per-country "Actual [2021]" (B) and "Actual [2022]" (D) service-unit
# figures, plus the handful of "Det. [2022]" (G) restatements that came in
# with this refresh. Row 6 (SES Area total) and the C/E/F/H/I ratio columns
# are formulas and recalculate automatically.
$ws.Range("B7").Value = 1799440
$ws.Range("D7").Value = 3247862
$ws.Range("B8").Value = 1166899
$ws.Range("D8").Value = 2096176
$ws.Range("G8").Value = 2107529
$ws.Range("B9").Value = 2269765
$ws.Range("D9").Value = 3870654
$ws.Range("B10").Value = 1518678
$ws.Range("D10").Value = 2228835
$ws.Range("G10").Value = 1581839
$ws.Range("B11").Value = 1266300
$ws.Range("D11").Value = 1788097
$ws.Range("G11").Value = 1837000
$ws.Range("B12").Value = 1280175
$ws.Range("B13").Value = 784993
$ws.Range("D13").Value = 1282410
$ws.Range("B14").Value = 466942
$ws.Range("D14").Value = 428511
$ws.Range("B15").Value = 494854
$ws.Range("D15").Value = 597862
$ws.Range("G15").Value = 894387
$ws.Range("B16").Value = 11180520
$ws.Range("D16").Value = 18897985
$ws.Range("B17").Value = 7678785
$ws.Range("D17").Value = 12390208
$ws.Range("B18").Value = 4048217
$ws.Range("D18").Value = 6416384
$ws.Range("G18").Value = 5861000
$ws.Range("B19").Value = 1726646
$ws.Range("D19").Value = 3184085
$ws.Range("B20").Value = 2419194
$ws.Range("D20").Value = 4233452
$ws.Range("G20").Value = 3991000
$ws.Range("B21").Value = 5782897
$ws.Range("D21").Value = 9561778
$ws.Range("B22").Value = 541944
$ws.Range("D22").Value = 465601
$ws.Range("G22").Value = 466000
$ws.Range("B23").Value = 443151
$ws.Range("D23").Value = 375999
$ws.Range("G23").Value = 372000
$ws.Range("B24").Value = 503699
$ws.Range("D24").Value = 666812
$ws.Range("G24").Value = 811000
$ws.Range("B25").Value = 1565320
$ws.Range("D25").Value = 2585835
$ws.Range("B26").Value = 1445483
$ws.Range("D26").Value = 2071287
$ws.Range("G26").Value = 2048218
$ws.Range("B27").Value = 2585928
$ws.Range("D27").Value = 3128964
$ws.Range("B28").Value = 1988333
$ws.Range("D28").Value = 3695099
$ws.Range("B29").Value = 2869907
$ws.Range("D29").Value = 4770304
$ws.Range("G29").Value = 4583000
$ws.Range("B30").Value = 611991
$ws.Range("D30").Value = 972528
$ws.Range("G30").Value = 798000
$ws.Range("B31").Value = 369971
$ws.Range("D31").Value = 595456
$ws.Range("G31").Value = 535978
$ws.Range("B32").Value = 1007563
$ws.Range("D32").Value = 1789655
$ws.Range("B33").Value = 6382913
$ws.Range("D33").Value = 11078709
$ws.Range("B34").Value = 1794889
$ws.Range("D34").Value = 2471898
$ws.Range("G34").Value = 2724000
$ws.Range("B35").Value = 897288
$ws.Range("D35").Value = 1544718
$ws.Range("G35").Value = 1594000
